$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.70%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.94"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.34%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.068"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.48%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05692"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.46%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.484"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.12%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8196"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.73%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8433"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.25%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.67%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06912"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.71%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02849"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.11%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09390"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.03%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001516"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.26%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04095"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-11.93%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.01008"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1,594.23%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006092"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.48%"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.509"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.52%"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.006"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.12%"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.316"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "12.68%"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3177"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.94%"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03195"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.30%"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1297"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.57%"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "MCDex"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.574"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.52%"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1373"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.69%"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001218"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.41%"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "HotbitToken"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.003965"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-13.57%"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "NitroEx"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009796"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2.00%"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "UpBots"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001937"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-0.08%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03701"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.01%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005798"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-6.07%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1055"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.13%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002299"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-8.08%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009400"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.49%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005198"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.90%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.08%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1199"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.08%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002473"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.66%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.08%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.08%"
